# Natmi following Dr Hou advice
# Rebuilds the Bmp6/Bmpr2 LR-pair table: adds "ECs" as a third sending
# cluster (alongside FAPs and sCs) and refreshes every computed metric,
# expanding the sheet from 8 data rows (A2:T9) to 12 data rows (A2:T13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Bmp6", "Bmpr2", "ECs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 33.66677566666667, 101.000327, 0.26047027611676, 0.26047027611676, 664.0605971936305, 5976.545374742675, 0.1519805918866407, 0.1519805918866408),
  @("ECs", "Bmp6", "Bmpr2", "FAPs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 43.683024, 131.049072, 0.3379631431162116, 0.3379631431162116, 861.6261709131999, 7754.635538218799, 0.1971965450047998, 0.1971965450047998),
  @("ECs", "Bmp6", "Bmpr2", "M2", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 17.657557, 52.97267100000001, 0.1366115006935798, 0.1366115006935798, 348.2866301928083, 3134.579671735275, 0.0797108101679343, 0.07971081016793431),
  @("ECs", "Bmp6", "Bmpr2", "sCs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 34.24645366666667, 102.739361, 0.2649550800734486, 0.2649550800734485, 675.4944607352805, 6079.450146617525, 0.154597409321608, 0.154597409321608),
  @("FAPs", "Bmp6", "Bmpr2", "ECs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 33.66677566666667, 101.000327, 0.26047027611676, 0.26047027611676, 25.870526758761, 232.834740828849, 0.005920872260501269, 0.00592087226050127),
  @("FAPs", "Bmp6", "Bmpr2", "FAPs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 43.683024, 131.049072, 0.3379631431162116, 0.3379631431162116, 33.56730244929599, 302.1057220436639, 0.007682399039848985, 0.007682399039848986),
  @("FAPs", "Bmp6", "Bmpr2", "M2", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 17.657557, 52.97267100000001, 0.1366115006935798, 0.1366115006935798, 13.568578867953, 122.117209811577, 0.003105380226031942, 0.003105380226031942),
  @("FAPs", "Bmp6", "Bmpr2", "sCs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 34.24645366666667, 102.739361, 0.2649550800734486, 0.2649550800734485, 26.315968144623, 236.843713301607, 0.006022818446979146, 0.006022818446979146),
  @("sCs", "Bmp6", "Bmpr2", "ECs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 33.66677566666667, 101.000327, 0.26047027611676, 0.26047027611676, 448.1618717526055, 4033.456845773449, 0.102568811969618, 0.102568811969618),
  @("sCs", "Bmp6", "Bmpr2", "FAPs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 43.683024, 131.049072, 0.3379631431162116, 0.3379631431162116, 581.495121287696, 5233.456091589264, 0.1330841990715628, 0.1330841990715628),
  @("sCs", "Bmp6", "Bmpr2", "M2", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 17.657557, 52.97267100000001, 0.1366115006935798, 0.1366115006935798, 235.0520250008197, 2115.468225007377, 0.05379531029961358, 0.05379531029961358),
  @("sCs", "Bmp6", "Bmpr2", "sCs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 34.24645366666667, 102.739361, 0.2649550800734486, 0.2649550800734485, 455.8783688732674, 4102.905319859407, 0.1043348523048614, 0.1043348523048614)
)

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
